# Apply trade #118 close (2026-02-17 16:03:42) to the live trading workbook.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Summary sheet - roll up numbers after the new closed trade.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1198.84   # Current Capital
$summary.Range("B4").Value = -1.17     # Total P&L $
$summary.Range("B6").Value = 118       # Total Trades
$summary.Range("B7").Value = 43        # Winning Trades
$summary.Range("B9").Value = 36.44     # Win Rate %

# ---------------------------------------------------------------------------
# 2. Strategy Status sheet - MarketMaking row (row 4).
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 98.84      # Capital
$status.Range("D4").Value = 118        # Trades
$status.Range("E4").Value = -1.17      # P&L $
$status.Range("F4").Value = -1.16      # P&L %
$status.Range("G4").Value = 36.44      # Win Rate %

# ---------------------------------------------------------------------------
# 3. Append the new trade row (#118) to both "All Trades" and "MarketMaking"
#    sheets, which mirror one another.
# ---------------------------------------------------------------------------
$newRow = 119

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item($newRow, 1).Value = 118      # A - Trade #

    $ws.Cells.Item($newRow, 2).NumberFormat = "@"
    $ws.Cells.Item($newRow, 2).Value = "2026-02-17"   # B - Date

    $ws.Cells.Item($newRow, 3).NumberFormat = "@"
    $ws.Cells.Item($newRow, 3).Value = "16:03:35"     # C - Time

    $ws.Cells.Item($newRow, 4).NumberFormat = "@"
    $ws.Cells.Item($newRow, 4).Value = "MarketMaking" # D - Strategy

    $ws.Cells.Item($newRow, 5).NumberFormat = "@"
    $ws.Cells.Item($newRow, 5).Value = "DOWN"         # E - Side

    $ws.Cells.Item($newRow, 6).Value = 0.079232       # F - Entry Price
    $ws.Cells.Item($newRow, 7).Value = 0.1            # G - Exit Price

    $ws.Cells.Item($newRow, 8).NumberFormat = "@"
    $ws.Cells.Item($newRow, 8).Value = "CLOSED"       # H - Status

    $ws.Cells.Item($newRow, 9).Value = 26.2114        # I - P&L %
    $ws.Cells.Item($newRow, 10).Value = 0.02          # J - P&L $
    $ws.Cells.Item($newRow, 11).Value = 98.84         # K - Capital After
    $ws.Cells.Item($newRow, 12).Value = 0             # L - Entry Slippage (bps)
    $ws.Cells.Item($newRow, 13).Value = 0             # M - Exit Slippage (bps)
    $ws.Cells.Item($newRow, 14).Value = 0.6           # N - Confidence

    $ws.Cells.Item($newRow, 15).NumberFormat = "@"
    $ws.Cells.Item($newRow, 15).Value = "Normal spread capture: 19600 bps" # O - Entry Reason

    $ws.Cells.Item($newRow, 16).NumberFormat = "@"
    $ws.Cells.Item($newRow, 16).Value = "early_exit"  # P - Exit Reason

    $ws.Cells.Item($newRow, 17).Value = 0.15          # Q - Duration (min)
}
